$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the business list values
$ws.Range("A2").Value = "jonas "
$ws.Range("A3").Value = "Cant reDd"

# Update selection to A3
$ws.Range("A3").Select()
